# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Note: several "Price" column values look like plain numbers (e.g. "235.65")
# but must stay as literal text (matching the source sheet's inlineStr cells),
# so a leading apostrophe forces text entry and the style is reset back to
# "Normal" afterwards so no stray number-format is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.577.19"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "2.376.95"
$ws.Range("E3").Value = "  +4.91%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'235.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").Value = "'0.656"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.10%  "
$ws.Range("E7").Value = "  +12.87%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.470"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.08%  "
$ws.Range("D10").Value = "'0.0974"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").Value = "'57.00"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.04%  "
$ws.Range("D12").Value = "'27.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("D13").Value = "2.741.63"
$ws.Range("E13").Value = "  +5.52%  "
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").Value = "'15.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "'6.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").Value = "'0.856"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("D18").Value = "2.385.58"
$ws.Range("E18").Value = "  +5.17%  "
$ws.Range("D19").Value = "43.498.19"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("D20").Value = "0.0₃0995"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").Value = "'6.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.31%  "
$ws.Range("D22").Value = "'74.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").Value = "'251.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").Value = "'3.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +16.11%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  +2.30%  "
$ws.Range("E27").Value = "  +2.23%  "
$ws.Range("D28").Value = "'22.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.81%  "
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").Value = "'174.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").Value = "'1.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.88%  "
$ws.Range("E32").Value = "  -7.99%  "
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").Value = "'5.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("D36").Value = "'5.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("D37").Value = "'6.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("D38").Value = "'2.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.69%  "
$ws.Range("D39").Value = "'3.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").Value = "'0.0256"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "'19.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.39%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "'8.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("E44").Value = "  +9.06%  "
$ws.Range("D45").Value = "'4.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.04%  "
$ws.Range("D46").Value = "'99.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("E47").Value = "  +1.95%  "
$ws.Range("D48").Value = "'0.0954"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "1.450.32"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "2.606.71"
$ws.Range("E50").Value = "  +5.36%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "'2.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "
